$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "org_day0_fio2"
$ws.Cells.Item(2, 2).Value = 0.1491249491283778
$ws.Cells.Item(3, 1).Value = "org_day0_map"
$ws.Cells.Item(3, 2).Value = 0.0841673494482206
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "1"
$ws.Cells.Item(4, 2).Value = 0.08115907950006938
$ws.Cells.Item(5, 1).Value = "ord_day0_gcs"
$ws.Cells.Item(5, 2).Value = 0.06848292355351977
$ws.Cells.Item(6, 1).Value = "demo_ethnicity___8"
$ws.Cells.Item(6, 2).Value = 0.04891827761822882
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "3"
$ws.Cells.Item(7, 2).Value = 0.04673174918215717
$ws.Cells.Item(8, 1).Value = "bl_sao2"
$ws.Cells.Item(8, 2).Value = 0.03705808193311192
$ws.Cells.Item(9, 1).Value = "demo_ethnicity___2"
$ws.Cells.Item(9, 2).Value = 0.0354166697883726
$ws.Cells.Item(10, 1).Value = "bl_lab_ast"
$ws.Cells.Item(10, 2).Value = 0.03372205296026062
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "2"
$ws.Cells.Item(11, 2).Value = 0.03012874185734349
$ws.Cells.Item(12, 1).Value = "org_day0_platlet"
$ws.Cells.Item(12, 2).Value = 0.02877700121690686
$ws.Cells.Item(13, 1).Value = "bl_lab_inr"
$ws.Cells.Item(13, 2).Value = 0.02660267123976017
$ws.Cells.Item(14, 1).Value = "bl_resp_rate"
$ws.Cells.Item(14, 2).Value = 0.02472475861076671
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "0"
$ws.Cells.Item(15, 2).Value = 0.02464591553443264
$ws.Cells.Item(16, 1).Value = "org_day0_bilirubin"
$ws.Cells.Item(16, 2).Value = 0.02296282386347054
$ws.Cells.Item(17, 1).Value = "bl_lab_wbc"
$ws.Cells.Item(17, 2).Value = 0.0208669386002578
$ws.Cells.Item(18, 1).Value = "bl_kg"
$ws.Cells.Item(18, 2).Value = 0.02079425108012155
$ws.Cells.Item(19, 1).Value = "demo_ethnicity___6"
$ws.Cells.Item(19, 2).Value = 0.02058454189865848
$ws.Cells.Item(20, 1).Value = "bl_lab_creatinine"
$ws.Cells.Item(20, 2).Value = 0.02015764660307344
$ws.Cells.Item(21, 1).Value = "demo_age_years"
$ws.Cells.Item(21, 2).Value = 0.02002214287272198
$ws.Cells.Item(22, 1).Value = "bl_lab_alt"
$ws.Cells.Item(22, 2).Value = 0.01821716182995778
$ws.Cells.Item(23, 1).Value = "bl_temp"
$ws.Cells.Item(23, 2).Value = 0.01785529013747142
$ws.Cells.Item(24, 1).Value = "bl_lab_troponin"
$ws.Cells.Item(24, 2).Value = 0.01642353345176189
$ws.Cells.Item(25, 1).Value = "bl_lab_haemo"
$ws.Cells.Item(25, 2).Value = 0.0163873361309962
$ws.Cells.Item(26, 1).Value = "bl_hr"
$ws.Cells.Item(26, 2).Value = 0.01607783475662551
$ws.Cells.Item(27, 1).Value = "co_smoking"
$ws.Cells.Item(27, 2).Value = 0.01440494081565701
$ws.Cells.Item(28, 1).Value = "demo_ethnicity___7"
$ws.Cells.Item(28, 2).Value = 0.01358411040334863
$ws.Cells.Item(29, 1).Value = "demo_ethnicity___3"
$ws.Cells.Item(29, 2).Value = 0.01309715945031404
$ws.Cells.Item(30, 1).Value = "demo_ethnicity___5"
$ws.Cells.Item(30, 2).Value = 0.009477893316570627
$ws.Cells.Item(31, 1).Value = "bl_oxy_status"
$ws.Cells.Item(31, 2).Value = 0.00675122437196147
$ws.Cells.Item(32, 1).Value = "demo_ethnicity___4"
$ws.Cells.Item(32, 2).Value = 0.004967025990691508
$ws.Cells.Item(33, 1).Value = "co_other___1"
$ws.Cells.Item(33, 2).Value = 0.004788206443707443
$ws.Cells.Item(34, 1).Value = "co_dementia___1"
$ws.Cells.Item(34, 2).Value = 0.002919716411104313

# Row 35 no longer holds data (sheet shrank from 35 to 34 data rows)
$ws.Rows.Item(35).Delete()
